# This script updates the "cryptos" worksheet with refreshed price / volume
# data scraped on Sat Nov 16 06:39:25 UTC 2024, matching the upstream
# GitHub Actions commit.
#
# Several coins (WrappedBTC/Toncoin, InternetComputer/Binance-Peg BSC-USD/
# Bittensor, Stellar/Hedera) also swapped ranking order, so their whole
# rows (Coin name, Link, Price, Volume) are rewritten in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain text value (coin name / URL) - these never look like
# numbers so a normal assignment is safe and keeps cell formatting untouched.
function Set-PlainValue([string]$cellRef, [string]$text) {
    $ws.Range($cellRef).Value = $text
}

# Helper: write a value that must stay plain text even though it may look
# like a number (e.g. "91.276.23" or "5.38"). Forcing the number format to
# "@" (Text) before assignment stops Excel's automatic type conversion from
# turning it into a numeric value, and resetting the style back to "Normal"
# afterwards keeps the cell's formatting identical to the original
# (unstyled) cell instead of leaving a stray text-format style behind.
function Set-TextValue([string]$cellRef, [string]$text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "91.276.23"
Set-TextValue "E2" "  +3.95%  "
Set-TextValue "D3" "3.119.75"
Set-TextValue "E3" "  +2.18%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "219.19"
Set-TextValue "E5" "  +4.78%  "
Set-TextValue "D6" "623.46"
Set-TextValue "E6" "  +1.05%  "
Set-TextValue "D7" "0.388"
Set-TextValue "E7" "  +5.39%  "
Set-TextValue "D8" "0.954"
Set-TextValue "E8" "  +17.98%  "
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "3.117.76"
Set-TextValue "E10" "  +2.32%  "
Set-TextValue "D11" "0.723"
Set-TextValue "E11" "  +20.25%  "
Set-TextValue "E12" "  +5.86%  "
Set-TextValue "D13" "0.0000255"
Set-TextValue "E13" "  +7.84%  "
Set-TextValue "D14" "34.45"
Set-TextValue "E14" "  +8.36%  "
Set-PlainValue "B15" "WrappedBTC"
Set-PlainValue "C15" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D15" "91.097.61"
Set-TextValue "E15" "  +3.66%  "
Set-PlainValue "B16" "Toncoin"
Set-PlainValue "C16" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D16" "5.38"
Set-TextValue "E16" "  +2.27%  "
Set-TextValue "D17" "3.698.15"
Set-TextValue "E17" "  +1.94%  "
Set-TextValue "D18" "3.133.23"
Set-TextValue "E18" "  +2.14%  "
Set-TextValue "D19" "3.77"
Set-TextValue "E19" "  +16.59%  "
Set-TextValue "E20" "  +10.31%  "
Set-TextValue "D21" "14.06"
Set-TextValue "E21" "  +6.73%  "
Set-TextValue "D22" "435.43"
Set-TextValue "E22" "  +4.06%  "
Set-TextValue "D23" "8.73"
Set-TextValue "E23" "  +7.75%  "
Set-TextValue "D24" "5.18"
Set-TextValue "E24" "  +6.48%  "
Set-TextValue "D25" "6.10"
Set-TextValue "E25" "  +11.80%  "
Set-TextValue "D26" "86.07"
Set-TextValue "E26" "  +4.93%  "
Set-TextValue "D27" "12.14"
Set-TextValue "E27" "  +3.82%  "
Set-TextValue "D28" "3.288.93"
Set-TextValue "E28" "  +1.88%  "
Set-TextValue "E29" "  +0.12%  "
Set-TextValue "D30" "0.167"
Set-TextValue "E30" "  +0.19%  "
Set-PlainValue "B31" "Binance-PegBSC-USD"
Set-PlainValue "C31" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -6.55%  "
Set-PlainValue "B32" "InternetComputer(DFINITY)"
Set-PlainValue "C32" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D32" "8.95"
Set-TextValue "E32" "  +11.74%  "
Set-PlainValue "B33" "Bittensor"
Set-PlainValue "C33" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D33" "528.43"
Set-TextValue "E33" "  +4.28%  "
Set-TextValue "E34" "  +6.84%  "
Set-TextValue "D35" "7.10"
Set-TextValue "E35" "  +6.50%  "
Set-TextValue "D36" "0.142"
Set-TextValue "E36" "  +8.05%  "
Set-TextValue "D37" "23.52"
Set-TextValue "E37" "  +6.05%  "
Set-TextValue "D38" "1.85"
Set-TextValue "E38" "  +3.51%  "
Set-TextValue "E39" "  +3.32%  "
Set-TextValue "E40" "  +0.36%  "
Set-TextValue "E41" "  +0.03%  "
Set-PlainValue "B42" "Hedera"
Set-PlainValue "C42" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D42" "0.0832"
Set-TextValue "E42" "  +24.28%  "
Set-PlainValue "B43" "Stellar"
Set-PlainValue "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D43" "0.149"
Set-TextValue "E43" "  +12.20%  "
Set-TextValue "D45" "0.377"
Set-TextValue "E45" "  +5.11%  "
Set-TextValue "D46" "1.90"
Set-TextValue "E46" "  +6.23%  "
Set-TextValue "D47" "147.04"
Set-TextValue "E47" "  -0.12%  "
Set-TextValue "D48" "44.01"
Set-TextValue "E48" "  +1.58%  "
Set-TextValue "D49" "1.29"
Set-TextValue "E49" "  +10.12%  "
Set-TextValue "D50" "166.22"
Set-TextValue "E50" "  +7.26%  "
Set-TextValue "D51" "4.17"
Set-TextValue "E51" "  +6.99%  "
